# Update recognition log: refresh attendee names/timestamps for existing rows (2-22)
# and append new recognition entries (rows 23-53) captured on 2024-12-03.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'DinhNhatKy'
$ws.Range("B2").Value = '2024-12-03 16:31:08'
# Row 3
$ws.Range("A3").Value = 'DinhNhatKy'
$ws.Range("B3").Value = '2024-12-03 16:31:25'
# Row 4
$ws.Range("A4").Value = 'DinhNhatKy'
$ws.Range("B4").Value = '2024-12-03 16:31:38'
$ws.Range("C4").Value = '11/15 (73.33%)'
# Row 5
$ws.Range("A5").Value = 'DinhNhatKy'
$ws.Range("B5").Value = '2024-12-03 16:31:54'
# Row 6
$ws.Range("B6").Value = '2024-12-03 16:32:13'
$ws.Range("C6").Value = '6/15 (40.00%)'
# Row 7
$ws.Range("B7").Value = '2024-12-03 16:32:24'
$ws.Range("C7").Value = '2/15 (13.33%)'
# Row 8
$ws.Range("A8").Value = 'Unknown'
$ws.Range("B8").Value = '2024-12-03 16:33:06'
$ws.Range("C8").Value = '0/15 (0.00%)'
# Row 9
$ws.Range("A9").Value = 'Quang lê'
$ws.Range("B9").Value = '2024-12-03 16:33:16'
# Row 10
$ws.Range("A10").Value = 'Quang lê'
$ws.Range("B10").Value = '2024-12-03 16:33:33'
$ws.Range("C10").Value = '14/15 (93.33%)'
# Row 11
$ws.Range("A11").Value = 'Quang lê'
$ws.Range("B11").Value = '2024-12-03 16:33:57'
$ws.Range("C11").Value = '15/15 (100.00%)'
# Row 12
$ws.Range("A12").Value = 'Unknown'
$ws.Range("B12").Value = '2024-12-03 16:34:16'
$ws.Range("C12").Value = '12/15 (80.00%)'
# Row 13
$ws.Range("B13").Value = '2024-12-03 16:34:30'
$ws.Range("C13").Value = '15/15 (100.00%)'
# Row 14
$ws.Range("A14").Value = 'Unknown'
$ws.Range("B14").Value = '2024-12-03 16:34:40'
$ws.Range("C14").Value = '13/15 (86.67%)'
# Row 15
$ws.Range("A15").Value = 'Unknown'
$ws.Range("B15").Value = '2024-12-03 16:34:53'
$ws.Range("C15").Value = '14/15 (93.33%)'
# Row 16
$ws.Range("A16").Value = 'Unknown'
$ws.Range("B16").Value = '2024-12-03 16:35:11'
$ws.Range("C16").Value = '6/15 (40.00%)'
# Row 17
$ws.Range("A17").Value = 'DinhNhatKy'
$ws.Range("B17").Value = '2024-12-03 16:35:32'
$ws.Range("C17").Value = '15/15 (100.00%)'
# Row 18
$ws.Range("B18").Value = '2024-12-03 16:35:44'
# Row 19
$ws.Range("A19").Value = 'DinhNhatKy'
$ws.Range("B19").Value = '2024-12-03 16:35:54'
# Row 20
$ws.Range("A20").Value = 'DinhNhatKy'
$ws.Range("B20").Value = '2024-12-03 16:36:05'
# Row 21
$ws.Range("B21").Value = '2024-12-03 16:37:24'
$ws.Range("C21").Value = '15/15 (100.00%)'
# Row 22
$ws.Range("A22").Value = 'DinhNhatKy'
$ws.Range("B22").Value = '2024-12-03 16:37:38'
$ws.Range("C22").Value = '13/15 (86.67%)'
# Row 23
$ws.Range("A23").Value = 'Anh Hung'
$ws.Range("B23").Value = '2024-12-03 16:38:04'
$ws.Range("C23").Value = '15/15 (100.00%)'
# Row 24
$ws.Range("A24").Value = 'Anh Hung'
$ws.Range("B24").Value = '2024-12-03 16:38:44'
$ws.Range("C24").Value = '15/15 (100.00%)'
# Row 25
$ws.Range("A25").Value = 'chipu'
$ws.Range("B25").Value = '2024-12-03 16:40:06'
$ws.Range("C25").Value = '15/15 (100.00%)'
# Row 26
$ws.Range("A26").Value = 'Unknown'
$ws.Range("B26").Value = '2024-12-03 16:41:18'
$ws.Range("C26").Value = '15/15 (100.00%)'
# Row 27
$ws.Range("A27").Value = 'Unknown'
$ws.Range("B27").Value = '2024-12-03 16:41:28'
$ws.Range("C27").Value = '15/15 (100.00%)'
# Row 28
$ws.Range("A28").Value = 'Unknown'
$ws.Range("B28").Value = '2024-12-03 16:41:42'
$ws.Range("C28").Value = '15/15 (100.00%)'
# Row 29
$ws.Range("A29").Value = 'DinhNhatKy'
$ws.Range("B29").Value = '2024-12-03 16:47:13'
$ws.Range("C29").Value = '15/15 (100.00%)'
# Row 30
$ws.Range("A30").Value = 'DinhNhatKy'
$ws.Range("B30").Value = '2024-12-03 16:47:25'
$ws.Range("C30").Value = '15/15 (100.00%)'
# Row 31
$ws.Range("A31").Value = 'DinhNhatKy'
$ws.Range("B31").Value = '2024-12-03 16:47:37'
$ws.Range("C31").Value = '14/15 (93.33%)'
# Row 32
$ws.Range("A32").Value = 'Unknown'
$ws.Range("B32").Value = '2024-12-03 16:48:10'
$ws.Range("C32").Value = '10/15 (66.67%)'
# Row 33
$ws.Range("A33").Value = 'Vanh Leg'
$ws.Range("B33").Value = '2024-12-03 16:48:23'
$ws.Range("C33").Value = '13/15 (86.67%)'
# Row 34
$ws.Range("A34").Value = 'Unknown'
$ws.Range("B34").Value = '2024-12-03 16:48:42'
$ws.Range("C34").Value = '9/15 (60.00%)'
# Row 35
$ws.Range("A35").Value = 'Vanh Leg'
$ws.Range("B35").Value = '2024-12-03 16:49:01'
$ws.Range("C35").Value = '13/15 (86.67%)'
# Row 36
$ws.Range("A36").Value = 'DinhNhatKy'
$ws.Range("B36").Value = '2024-12-03 16:49:34'
$ws.Range("C36").Value = '15/15 (100.00%)'
# Row 37
$ws.Range("A37").Value = 'Ta Hoang Giang'
$ws.Range("B37").Value = '2024-12-03 16:50:10'
$ws.Range("C37").Value = '15/15 (100.00%)'
# Row 38
$ws.Range("A38").Value = 'Nguyen Van Tinh'
$ws.Range("B38").Value = '2024-12-03 16:50:40'
$ws.Range("C38").Value = '15/15 (100.00%)'
# Row 39
$ws.Range("A39").Value = 'DinhNhatKy'
$ws.Range("B39").Value = '2024-12-03 16:51:01'
$ws.Range("C39").Value = '13/15 (86.67%)'
# Row 40
$ws.Range("A40").Value = 'Unknown'
$ws.Range("B40").Value = '2024-12-03 16:51:25'
$ws.Range("C40").Value = '0/15 (0.00%)'
# Row 41
$ws.Range("A41").Value = 'Nguyen Van Tinh'
$ws.Range("B41").Value = '2024-12-03 16:51:38'
$ws.Range("C41").Value = '15/15 (100.00%)'
# Row 42
$ws.Range("A42").Value = 'Nguyen Van Tinh'
$ws.Range("B42").Value = '2024-12-03 16:52:24'
$ws.Range("C42").Value = '13/15 (86.67%)'
# Row 43
$ws.Range("A43").Value = 'NguyenHuuDuc'
$ws.Range("B43").Value = '2024-12-03 16:52:25'
$ws.Range("C43").Value = '15/15 (100.00%)'
# Row 44
$ws.Range("A44").Value = 'DinhNhatKy'
$ws.Range("B44").Value = '2024-12-03 16:52:45'
$ws.Range("C44").Value = '15/15 (100.00%)'
# Row 45
$ws.Range("A45").Value = 'NguyenHuuDuc'
$ws.Range("B45").Value = '2024-12-03 16:53:02'
$ws.Range("C45").Value = '15/15 (100.00%)'
# Row 46
$ws.Range("A46").Value = 'Nguyen Van Tinh'
$ws.Range("B46").Value = '2024-12-03 16:53:25'
$ws.Range("C46").Value = '15/15 (100.00%)'
# Row 47
$ws.Range("A47").Value = 'sontung'
$ws.Range("B47").Value = '2024-12-03 16:55:26'
$ws.Range("C47").Value = '15/15 (100.00%)'
# Row 48
$ws.Range("A48").Value = 'sontung'
$ws.Range("B48").Value = '2024-12-03 16:55:39'
$ws.Range("C48").Value = '15/15 (100.00%)'
# Row 49
$ws.Range("A49").Value = 'sontung'
$ws.Range("B49").Value = '2024-12-03 16:55:53'
$ws.Range("C49").Value = '15/15 (100.00%)'
# Row 50
$ws.Range("A50").Value = 'sontung'
$ws.Range("B50").Value = '2024-12-03 16:56:12'
$ws.Range("C50").Value = '15/15 (100.00%)'
# Row 51
$ws.Range("A51").Value = 'Chi Dân'
$ws.Range("B51").Value = '2024-12-03 16:56:27'
$ws.Range("C51").Value = '15/15 (100.00%)'
# Row 52
$ws.Range("A52").Value = 'Karik'
$ws.Range("B52").Value = '2024-12-03 16:57:00'
$ws.Range("C52").Value = '15/15 (100.00%)'
# Row 53
$ws.Range("A53").Value = 'Karik'
$ws.Range("B53").Value = '2024-12-03 16:57:17'
$ws.Range("C53").Value = '15/15 (100.00%)'
